# EchoPlay Gantt chart update
# - Move the displayed week window forward (Anzeigewoche 4 -> 6), which
#   cascades through the sheet's date formulas (rows 4 & 5 headers).
# - Mark the "Ersatz Chip Plannen" task as fully complete (50% -> 100%).
# - Refresh the on-screen view (zoom + active cell) to match where the
#   author was last looking in the plan.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Projektplan")

# Anzeigewoche (display week) - drives the rolling 8-week window shown
# in the Gantt header via formulas that reference this cell.
$ws.Range("E4").Value = 6

# Task progress for "Ersatz Chip Plannen" (row 13) -> 100%.
$ws.Range("D13").Value = 1

# Restore focus/view state: zoomed out a bit and scrolled further down
# into the task list.
$ws.Activate()
[void]$ws.Range("L18").Select()
$excel.ActiveWindow.Zoom = 67
